$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (want-to-go count) column F
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F5").Value = 972
$wsExhibition.Range("F6").Value = 2322
$wsExhibition.Range("F7").Value = 200

# Sheet "全部类型" (All types) - same underlying rows duplicated, update column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 972
$wsAll.Range("F8").Value = 2322
$wsAll.Range("F10").Value = 200
